$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-12, columns D (Fecha), L (Calidad), M (Volumen),
# N (Precio minimo), O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg)
$rows = @(
    @{ Row = 2;  D = 44517; L = "Especial"; M = 100; N = 27000; O = 27000; P = 27000; S = 2700 },
    @{ Row = 3;  D = 44517; L = "Primera";  M = 30;  N = 25000; O = 25000; P = 25000; S = 2500 },
    @{ Row = 4;  D = 44432; L = "Primera";  M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 5;  D = 44473; L = "Primera";  M = 180; N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 6;  D = 44476; L = "Primera";  M = 120; N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 7;  D = 44466; L = "Primera";  M = 60;  N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 8;  D = 44434; L = "Primera";  M = 20;  N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 9;  D = 44511; L = "Primera";  M = 120; N = 28000; O = 28000; P = 28000; S = 2800 },
    @{ Row = 10; D = 44435; L = "Primera";  M = 40;  N = 20000; O = 20000; P = 20000; S = 2000 },
    @{ Row = 11; D = 44503; L = "Primera";  M = 60;  N = 30000; O = 30000; P = 30000; S = 3000 },
    @{ Row = 12; D = 44503; L = "Segunda";  M = 50;  N = 25000; O = 25000; P = 25000; S = 2500 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Calidad
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Volumen
    $ws.Cells.Item($row, 14).Value = $r.N   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value = $r.O   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $r.S   # S - Precio $/Kg
}
